$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four category-header rows (sexo / cor ou raça / grupos de idade /
# classes de rendimento mensal domiciliar per capita) and the trailing
# source-note row only ever held a label in column A with no data in B:H.
# The fix merges each label into the row immediately below it (which held
# the first data series of that category) by deleting the empty-data label
# row, which shifts everything below up by one. Delete from the bottom up
# so earlier row numbers stay stable while we work.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
